# Reproduce the NSantiam 2010-18 sheet edit: add a new "Baseline 2010-18 C267+"
# simulation row (row 4) plus a blank formatted spacer row (row 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2010-18")

# --- Row 4: new scenario results -----------------------------------------
$ws.Range("A4").Value = "CW3M"
$ws.Range("B4").Value = "Baseline 2010-18 C267+"
$ws.Range("C4").Value = "2010-18"

$ws.Range("D4").NumberFormat = "0.00"
$ws.Range("D4").Value = 677.97837322222222

$ws.Range("E4").NumberFormat = "0.00"
$ws.Range("E4").Value = 2094.2995878888887

$ws.Range("F4").NumberFormat = "0.00"
$ws.Range("F4").Interior.Color = 65535
$ws.Range("F4").Value = 4.820043222222222

$ws.Range("G4").NumberFormat = "0.00"
$ws.Range("G4").Interior.Color = 65535
$ws.Range("G4").Value = 232.21855144444442

$ws.Range("H4").NumberFormat = "0.00"
$ws.Range("H4").Value = 0

$ws.Range("I4").NumberFormat = "0.00"
$ws.Range("I4").Interior.Color = 65535
$ws.Range("I4").Value = 6.3389989999999994

$ws.Range("J4").NumberFormat = "0.00"
$ws.Range("J4").Value = 0

$ws.Range("K4").NumberFormat = "0.00"
$ws.Range("K4").Interior.Color = 65535
$ws.Range("K4").Value = 586.02156566666656

$ws.Range("L4").NumberFormat = "0.00"
$ws.Range("L4").Value = 96.631732222222212

$ws.Range("M4").NumberFormat = "0.00"
$ws.Range("M4").Interior.Color = 65535
$ws.Range("M4").Value = 1650.8734266666665

$ws.Range("N4").NumberFormat = "0.00"
$ws.Range("N4").Value = 682.41798233333327

$ws.Range("O4").NumberFormat = "0"
$ws.Range("O4").Interior.Color = 65535
$ws.Range("O4").Value = 12820.605631666667

$ws.Range("P4").NumberFormat = "0"
$ws.Range("P4").Value = 2216.7525497777779

$ws.Range("Q4").NumberFormat = "0.00"
$ws.Range("Q4").Value = 0.28915188888888882

# PowerShell literal parser here chokes on "E-6" scientific notation, so
# build the value arithmetically instead of typing it as a literal.
$ws.Range("R4").NumberFormat = "0.000000"
$ws.Range("R4").Value = (-2.4444444444444798 / 1000000)

# --- Row 5: blank formatted spacer row ------------------------------------
$ws.Range("D5:N5").NumberFormat = "0.00"
$ws.Range("O5").NumberFormat = "0"
$ws.Range("O5").Interior.Color = 65535
$ws.Range("P5").NumberFormat = "0"
$ws.Range("Q5").NumberFormat = "0.00"
$ws.Range("R5").NumberFormat = "0.000000"

# --- Selection as left by the author after editing ------------------------
$ws.Range("A6:XFD6").Select() | Out-Null
